$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows before the header row (old row 8 becomes row 15),
# pushing the header row (and everything after) down.
$ws.Rows("7:13").Insert()

# --- Filter-label column (A3:A13) ---
$ws.Range("A4").Value = "Charge Type"
$ws.Range("A5").Value = "Org"
$ws.Range("A6").Value = "Sector"
$ws.Range("A7").Value = "Subsector"
$ws.Range("A8").Value = "Division"
$ws.Range("A9").Value = "Section"
$ws.Range("A10").Value = "Budget Method"
$ws.Range("A11").Value = "Project-C"
$ws.Range("A12").Value = "Run By"
$ws.Range("A13").Value = "Run Date"

# A3:A11 use the wrapped "field label" look (teal fill + bold + wrap text);
# A12:A13 keep the plain (non-wrap) look the old "Run By"/"Run Date" rows
# had. Set explicitly on both groups so the formatting round-trips cleanly.
$ws.Range("A3:A11").WrapText = $true
$ws.Range("A12:A13").WrapText = $false

# --- Header row (now row 15) ---
$headers = @("Org","Division","Section Code","Section Name","Project C Code","Project C Name","Phase Code","Phase Name","Status","Budget Plan","Budget Release","PR Commit","PO Commit","Expense Commit","Total Commit","Actual","Total Spent","Balance","%Commitment","%Actual","%Usage")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(15, $col).Value = $headers[$i]
}
$ws.Range("A15:U15").WrapText = $true

$ws.Range("A16").Select()

# Mirrors the source file's trailing sentinel row at the very bottom of the
# sheet (row 1048576) -- a harmless no-content row pinning the sheet's row
# height, left over from the original template's save.
$ws.Rows(1048576).RowHeight = 12.8
